# Generate Report for Handback
# - Flip the "in sync" status to "not in sync" everywhere it is shown
#   (Overview!E/F for both rows, and the Status column (C) on the
#   zh-cn / de-de detail sheets - they all display the same text).
# - Record a fresh "Correspond Handback DateTime" for the second file
#   (c71c5377...) on both the zh-cn and de-de sheets.
# - Widen the "in/out of sync" status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns for both rows ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn / de-de detail sheets: Status column (C) for both rows ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- New handback timestamps for the c71c5377... file (row 3) ---
$wsZhCn.Range("K3").Value = "2016-10-21 04:38:58"
$wsDeDe.Range("K3").Value = "2016-10-21 04:39:16"

# --- Widen the status columns to fit the longer "not in sync" text ---
$wsOverview.Columns.Item(5).ColumnWidth = 32.6667
$wsOverview.Columns.Item(6).ColumnWidth = 32.6667
$wsZhCn.Columns.Item(3).ColumnWidth = 32.6667
$wsDeDe.Columns.Item(3).ColumnWidth = 32.6667
